# fix(Import): reliable detection of intended datetime type
#
# Inserts a new leading "Date and Time" column (full date+time values) and a
# trailing "time" column (time-only values), shifting the previously
# imported columns (Col1..truth/date) one position to the right, and fixes
# the boolean columns to use a plain General format instead of the
# LibreOffice "TRUE/FALSE" text-number-format trick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make room: insert a new column before column A. This shifts all
#    existing data (and its formatting) one column to the right.
$ws.Columns.Item(1).Insert()

# 2. New column A: "Date and Time" header + datetime values.
$ws.Range("A1").Value = "Date and Time"
$ws.Range("A1").Style = $ws.Range("B1").Style

$ws.Range("A2").Value = 44612.3625
$ws.Range("A3").Value = 42522.5673611111
$ws.Range("A2:A3").NumberFormat = "yyyy-mm-dd hh:mm"

# 3. New column J: "time" header + time-only values (fraction of a day).
$ws.Range("J1").Value = "time"
$ws.Range("J1").Style = $ws.Range("B1").Style

$ws.Range("J2").Value = 0.783333333333333
$ws.Range("J3").Value = 0.0576388888888889
$ws.Range("J2:J3").NumberFormat = "hh:mm"

# 4. Booleans (shifted from H to I) now get a plain General format instead
#    of the previous "TRUE";"TRUE";"FALSE" text trick, now that the cells
#    are properly typed as booleans.
$ws.Range("I2:I3").NumberFormat = "General"

# 5. Column A is a bit wider to fit the date+time values.
$ws.Columns.Item(1).ColumnWidth = 17.43

# 6. Restore a sane selection (matches the post-edit author's cursor spot).
$ws.Range("G7").Select()

Write-Host "applied datetime/time column fix"
